# edit.ps1 - applies the "Recolored navbar and quality of life stuff" change
# (the textual part of it: adding two bullets to slide 5's content list and
# rewording a bullet on slide 6) to the Caravanserai deck.

function Find-CharRange {
    param(
        $TextRange,
        [string]$Needle
    )
    $len = $TextRange.Length
    $n = $Needle.Length
    for ($pos = 1; $pos -le ($len - $n + 1); $pos++) {
        $cand = $TextRange.Characters($pos, $n)
        if ($cand.Text -eq $Needle) {
            return $cand
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 5 ("What we have Built"): add a blank bullet line and a new
# "Shared inventory between group members" bullet after the existing
# "Message and Barter system" line.
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$contentShape5 = $slide5.Shapes.Item("Content Placeholder 4")
$tr5 = $contentShape5.TextFrame.TextRange

$anchor5 = Find-CharRange $tr5 "Message and Barter system"
if ($null -eq $anchor5) {
    throw "Could not locate 'Message and Barter system' bullet on slide 5"
}

# InsertAfter on the shape's full text range appends new paragraphs that
# inherit the bullet/paragraph formatting of the last paragraph - exactly
# mirroring the blank-bullet-line pattern already used elsewhere in this
# placeholder.
$tr5.InsertAfter("`r`n`r`nShared inventory between group members")

# ---------------------------------------------------------------------
# Slide 6 ("What IS TO COME"): reword the "Quality of life and
# finalizing backend" bullet to "Finalize bartering process".
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$contentShape6 = $slide6.Shapes.Item("Content Placeholder 4")
$tr6 = $contentShape6.TextFrame.TextRange

$run1 = Find-CharRange $tr6 "Quality of life and "
if ($null -eq $run1) {
    throw "Could not locate 'Quality of life and ' run on slide 6"
}
$run1.Text = "Finalize "

# Recompute against a fresh TextRange since the text just shifted.
$tr6b = $contentShape6.TextFrame.TextRange
$run2 = Find-CharRange $tr6b "finalizing backend"
if ($null -eq $run2) {
    throw "Could not locate 'finalizing backend' run on slide 6"
}
$run2.Text = "bartering process"
